# Add batch sql test case about nacos (batch_009) to the mysql batchsql test sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Add Sub_component = "SingleTable" for the existing single-table rows (2-9).
$ws.Cells.Item(2,5).Value = "SingleTable"
$ws.Cells.Item(3,5).Value = "SingleTable"
$ws.Cells.Item(4,5).Value = "SingleTable"
$ws.Cells.Item(5,5).Value = "SingleTable"
$ws.Cells.Item(6,5).Value = "SingleTable"
$ws.Cells.Item(7,5).Value = "SingleTable"
$ws.Cells.Item(8,5).Value = "SingleTable"
$ws.Cells.Item(9,5).Value = "SingleTable"

# 2. Fix the expected-result csv paths for batch_007 / batch_008 (folder rename
#    from .../testdata/cases/... to .../testdata/mysqlcases/...).
$ws.Cells.Item(8,10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_007.csv"
$ws.Cells.Item(9,10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_008.csv"

# 3. Add a new row 10 for the batch_009 (nacos multi-table) test case.
$ws.Rows.Item(10).Insert() | Out-Null

$ws.Cells.Item(10,1).Value = "batch_009"
$ws.Cells.Item(10,2).Value = "y"
$ws.Cells.Item(10,3).Value = "批量操作语句9执行"
$ws.Cells.Item(10,4).Value = "batchsql"
$ws.Cells.Item(10,5).Value = "MultiTable"
$ws.Cells.Item(10,6).Value = "config_info,config_info_aggr,config_info_beta,config_info_tag,config_tags_relation,group_capacity,his_config_info,tenant_capacity,tenant_info,users,roles,permissions"
$ws.Cells.Item(10,8).Value = "batch_sql_09"
$ws.Cells.Item(10,9).Value = "select * from mysql_users"
$ws.Cells.Item(10,10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_009_01.csv"
$ws.Cells.Item(10,11).Value = "select * from mysql_roles"
$ws.Cells.Item(10,12).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/batchsql/expectedresult/batch_009_02.csv"
$ws.Cells.Item(10,13).Value = "csv_containsAll"

# 4. Add "justExec" to the Validation_type drop-down list (column M).
$dv = $ws.Range("M2:M1048576").Validation()
$dv.Modify(3, 1, 1, '"csv_equals,csv_containsAll,string_equals,effected_rows_assert,assertNull,justExec,SQLException"')

# 5. Restore the selection as it ended up after the edit.
$ws.Range("G19").Select() | Out-Null
